# Add 2022-Q3 data
# 1) "总计" (totals) sheet: insert a new row 2 for "2022-Q3" and shift the
#    existing quarters down, fixing up the running index in column A.
# 2) Insert a brand-new "2022-Q3" worksheet (positioned right after "总计",
#    before "2022-Q2") and fill it with the fund-holding detail rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Push existing data rows (old rows 2-7) down to rows 3-8, carrying their
# formatting (incl. the bold/bordered/centered style on column A) along.
$totals.Rows("2:2").Insert()

# The inserted row picked up a "data row" look-alike format on B:D; strip
# that so the new row matches the plain (unstyled) data cells elsewhere.
$totals.Range("B2:D2").ClearFormats()

# Give the new A2 the same bold/centered/bordered look as the other index
# cells in column A (copy formatting from the cell right below it).
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

# New 2022-Q3 summary row.
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 7
$totals.Range("D2").Value = 0.55

# Renumber the running index in column A for the rows that shifted down.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5
$totals.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet, right before "2022-Q2"
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($q2)
$ws.Name = "2022-Q3"

# Header row.
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Style the header row: bold, thin border all around, centered/top aligned
# - matches the look of the header row on the other quarter sheets.
$header = $ws.Range("B1:H1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Data rows.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'009715"
$ws.Range("C2").Value = "汇添富策略增长灵活配置混合"
$ws.Range("D2").Value = "'3.87"
$ws.Range("E2").Value = "'88.49"
$ws.Range("F2").Value = "'5.01"
$ws.Range("G2").Value = "'0.1939"
$ws.Range("H2").Value = 4

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'010676"
$ws.Range("C3").Value = "光大保德信新机遇混合"
$ws.Range("D3").Value = "'2.85"
$ws.Range("E3").Value = "'84.08"
$ws.Range("F3").Value = "'4.47"
$ws.Range("G3").Value = "'0.1274"
$ws.Range("H3").Value = 9

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'160218"
$ws.Range("C4").Value = "国泰国证房地产行业指数A"
$ws.Range("D4").Value = "'5.91"
$ws.Range("E4").Value = "'94.23"
$ws.Range("F4").Value = "'1.95"
$ws.Range("G4").Value = "'0.1152"
$ws.Range("H4").Value = 9

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'515060"
$ws.Range("C5").Value = "华夏中证全指房地产ETF"
$ws.Range("D5").Value = "'4.09"
$ws.Range("E5").Value = "'99.07"
$ws.Range("F5").Value = "'1.84"
$ws.Range("G5").Value = "'0.0753"
$ws.Range("H5").Value = 9

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'015042"
$ws.Range("C6").Value = "国泰国证房地产行业指数C"
$ws.Range("D6").Value = "'1.29"
$ws.Range("E6").Value = "'94.23"
$ws.Range("F6").Value = "'1.95"
$ws.Range("G6").Value = "'0.0252"
$ws.Range("H6").Value = 9

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'519172"
$ws.Range("C7").Value = "浦银安盛睿智精选灵活配置混合A"
$ws.Range("D7").Value = "'0.22"
$ws.Range("E7").Value = "'79.23"
$ws.Range("F7").Value = "'3.03"
$ws.Range("G7").Value = "'0.0067"
$ws.Range("H7").Value = 4

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'519173"
$ws.Range("C8").Value = "浦银安盛睿智精选灵活配置混合C"
$ws.Range("D8").Value = "'0.19"
$ws.Range("E8").Value = "'79.23"
$ws.Range("F8").Value = "'3.03"
$ws.Range("G8").Value = "'0.0058"
$ws.Range("H8").Value = 4

# Style column A (the running index) to match the bold/bordered/centered
# look used for that column on the other quarter sheets.
$idxCol = $ws.Range("A2:A8")
$idxCol.Font.Bold = $true
$idxCol.Borders.LineStyle = 1
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
